# Fix category case sensitivity issue - normalize default category to match
# ML model output. Adds the two new "Update" records (Alice Brown / Marketing
# and Charlie Wilson / Sales) that were produced with the corrected category
# casing, and widens the Name / Email columns to fit the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 3: Alice Brown / Marketing
$ws.Range("A3").Value = "Alice Brown"
$ws.Range("B3").Value = "alice@example.com"
$ws.Range("C3").Value = "Marketing"
$ws.Range("E3").Value = "Update"

# New row 4: Charlie Wilson / Sales
$ws.Range("A4").Value = "Charlie Wilson"
$ws.Range("B4").Value = "charlie@example.com"
$ws.Range("C4").Value = "Sales"
$ws.Range("E4").Value = "Update"

# Widen columns A (Name) and B (Email) so the longer new values fit.
$ws.Columns.Item(1).ColumnWidth = 12.877291666666666
$ws.Columns.Item(2).ColumnWidth = 19.877291666666668
